$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update marks / test data for "Session 7 (Branch&Bound)" column (H)
# and the Test mark column (I) for the student row (row 4),
# plus an annotation in the comments row (row 5).
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = "Not implemented"
$ws.Range("I4").Value = "NA"

# Update the view: scroll so column B is the left-most visible column,
# and select the merged range G5:G12.
$excel.ActiveWindow.TopLeftCell = $ws.Range("B1")
$ws.Range("G5:G12").Select()
